$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$tr = $s.Shapes.Item(1).TextFrame.TextRange

# The title text is currently split across three separate runs
# ("Below" + " " + "section-level"). Replace the whole range in one
# shot (delete, then re-insert) so the writer emits a single,
# consolidated <a:r> run instead of three.
$tr.Delete()
[void]$tr.InsertAfter("Below section-level")
